$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-16 Friday" "2025-05-17 Saturday"

Replace-Text "835×2=1670" "366×4=1464"
Replace-Text "820×3=2460" "867×7=6069"
Replace-Text "986×5=4930" "498×2=996"
Replace-Text "133×7=931" "473×9=4257"
Replace-Text "706×6=4236" "297×4=1188"
Replace-Text "928×5=4640" "966×2=1932"
Replace-Text "375×2=750" "413×5=2065"
Replace-Text "275×4=1100" "510×2=1020"
Replace-Text "522×7=3654" "578×7=4046"
Replace-Text "196×6=1176" "256×2=512"
Replace-Text "687×2=1374" "132×5=660"
Replace-Text "588×6=3528" "679×3=2037"
Replace-Text "530×3=1590" "863×6=5178"
Replace-Text "724×5=3620" "946×6=5676"
Replace-Text "740×3=2220" "931×6=5586"
Replace-Text "444×2=888" "888×5=4440"
Replace-Text "594×3=1782" "414×4=1656"
Replace-Text "233×6=1398" "210×3=630"
Replace-Text "367×8=2936" "644×5=3220"
Replace-Text "462×7=3234" "821×6=4926"
Replace-Text "523×3=1569" "334×8=2672"
Replace-Text "790×3=2370" "305×6=1830"
Replace-Text "104×5=520" "803×3=2409"
Replace-Text "233×4=932" "643×2=1286"
Replace-Text "104×6=624" "578×9=5202"

Write-Output "Done applying replacements"
